# "update scripts wuth new tpm"
# The underlying analysis script was re-run with updated TPM values. This
# changed which Sending/Target cluster combinations passed the significance
# filter: all rows whose Target cluster was "ECs" dropped out, leaving only
# the four rows whose Target cluster is "FAPs" - and the numeric metrics for
# those surviving rows were recomputed against the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the rows with Target cluster = "ECs" (original rows 2, 4, 6, 8).
# Delete from the bottom up so earlier row numbers remain valid targets.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(2).Delete()

$columns = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Recomputed metrics (columns E:T) for the four remaining rows (now rows 2-5),
# in row order.
$newValues = @(
    @(3, 1, 3.560699333333333, 10.682098, 0.2516303646515017, 0.2516303646515017, 2, 0.6666666666666666, 0.6668756666666665, 2.000627, 1, 1, 2.374543741716222, 21.370893675446, 0.2516303646515017, 0.2516303646515017),
    @(3, 1, 3.046736666666666, 9.14021, 0.2153092375010323, 0.2153092375010323, 2, 0.6666666666666666, 0.6668756666666665, 2.000627, 1, 1, 2.031794545741111, 18.28615091167, 0.2153092375010323, 0.2153092375010323),
    @(3, 1, 4.835201333333333, 14.505604, 0.3416978971743455, 0.3416978971743456, 2, 0.6666666666666666, 0.6668756666666665, 2.000627, 1, 1, 3.224478112634221, 29.02030301370799, 0.3416978971743455, 0.3416978971743456),
    @(3, 1, 2.707878, 8.123634000000001, 0.1913625006731204, 0.1913625006731204, 2, 0.6666666666666666, 0.6668756666666665, 2.000627, 1, 1, 1.805817946502, 16.252361518518, 0.1913625006731204, 0.1913625006731204)
)

for ($r = 0; $r -lt $newValues.Length; $r++) {
    $row = $r + 2
    $values = $newValues[$r]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value2 = $values[$i]
    }
}
